$wb = $excel.ActiveWorkbook

# Helper: write a value that looks numeric (e.g. "1.1") into a cell while
# forcing it to stay a text/string cell (matching the workbook's existing
# convention of storing these expression-evaluation numbers as strings).
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# NOTE: worksheet names are looked up case-insensitively, and this workbook
# has two sheets whose names differ only by case ("Vector_bf" / "Vector_BF"),
# so sheets are addressed by their (1-based) tab position instead of by name
# to avoid ambiguity:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Restricciones_del_lider ---
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2.Range("A2") "-2.1 + x"
Set-TextValue $ws2.Range("B2") "1.1"
Set-TextValue $ws2.Range("D2") "0.74"
Set-TextValue $ws2.Range("A3") "2.1 - x"
Set-TextValue $ws2.Range("B3") "-3.1"
Set-TextValue $ws2.Range("D3") "0.27"

# --- Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3.Range("A2") "-51.117647058823536 + 15.490196078431374y"
Set-TextValue $ws3.Range("B2") "50.117647058823536"
Set-TextValue $ws3.Range("D2") "0.22"
Set-TextValue $ws3.Range("F2") "7.9"
Set-TextValue $ws3.Range("A3") "1.617 - 0.49y"
Set-TextValue $ws3.Range("B3") "-2.617"
Set-TextValue $ws3.Range("D3") "0.66"
Set-TextValue $ws3.Range("F3") "-7.199999999999999"

# --- Punto_modificado ---
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "2.1"
Set-TextValue $ws4.Range("B2") "3.3000000000000003"

# --- Vector_bf ---
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-11.874443137254906"

# --- Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-3.66"
Set-TextValue $ws6.Range("A3") "-11.191"

# --- Vector_Alpha (this one stays a genuine number, not text) ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 0.51
